# Weekly update: a new price observation (row) is inserted at row 168,
# pushing the existing rows 168-254 down to 169-255 (dimension grows from
# A1:R254 to A1:R255).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 168.
$ws.Rows.Item(168).Insert()

# Populate the newly inserted row 168 with the new weekly record.
$ws.Cells.Item(168, 1).Value  = 7
$ws.Cells.Item(168, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(168, 3).Value  = "Ñuble"
$ws.Cells.Item(168, 4).Value  = 44460
$ws.Cells.Item(168, 5).Value  = 16
$ws.Cells.Item(168, 6).Value  = 100114014
$ws.Cells.Item(168, 7).Value  = "Betarraga"
$ws.Cells.Item(168, 8).Value  = "Sin especificar"
$ws.Cells.Item(168, 9).Value  = "Primera"
$ws.Cells.Item(168, 10).Value = 300
$ws.Cells.Item(168, 11).Value = 750
$ws.Cells.Item(168, 12).Value = 800
$ws.Cells.Item(168, 13).Value = 775
$ws.Cells.Item(168, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(168, 15).Value = "Región del Maule"
$ws.Cells.Item(168, 16).Value = 155
$ws.Cells.Item(168, 17).Value = 5
$ws.Cells.Item(168, 18).Value = "Hortaliza"

# Match the date formatting style used by the rest of column D.
$ws.Cells.Item(168, 4).NumberFormat = $ws.Cells.Item(169, 4).NumberFormat
